$d = $word.ActiveDocument

# --- Locate the run of text to split -------------------------------------
# "Clone the repository using git clone " needs to become three runs
# (same run formatting throughout):
#   1) "Clone the repository using "
#   2) "the command $ "
#   3) "git clone "
# with the document's "_GoBack" bookmark relocated to the boundary between
# runs 2 and 3.

$rng = $d.Content
$found = $rng.Find.Execute(
    "Clone the repository using git clone ",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find target sentence to edit"
}

# Shrink the matched range down to just the first chunk of text; this
# replaces the whole original run with the shorter text (single run).
$rng.Text = "Clone the repository using "

# Append the second chunk right after it.
$rng.Collapse(0)
$boundary1 = $rng.End
$rng.InsertAfter("the command `$ ")

# Append the third chunk right after that.
$rng.Collapse(0)
$boundary2 = $rng.End
$rng.InsertAfter("git clone ")

# --- Force a run split between chunk 1 and chunk 2 ------------------------
# Adjacent runs sharing identical formatting get coalesced on save, so drop
# a throwaway bookmark on the boundary and remove it again; the split it
# creates survives the bookmark's own removal.
$tempRange = $d.Range($boundary1, $boundary1)
$d.Bookmarks.Add("TempSplitMark", $tempRange) | Out-Null
$d.Bookmarks("TempSplitMark").Delete()

# --- Relocate "_GoBack" to the boundary between chunk 2 and chunk 3 -------
# (Word automatically drops the bookmark from its old location - the end of
# the "Updated README.md file" paragraph - when it's re-added here.)
$goBackRange = $d.Range($boundary2, $boundary2)
$d.Bookmarks.Add("_GoBack", $goBackRange) | Out-Null
